# conditions.xlsx update
#
# Commit: "Add: better quality pictures. Update: scales. TODO: write out the
# scales to a file. better quality radio buttons and scales in general (add
# question about price to carbon questions)."
#
# The "better quality pictures" part of the commit renames the smartband
# product image from the lower quality "smartband2.png" to the higher
# quality "smartband.png". That filename is stored as the "image" column
# (column E) value for the "smartband" product row (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 is the smartband product row; column E holds the image path string.
$ws.Range("E5").Value = "png/smartband.png"

# The author's cursor ended up one row below the edited cell afterwards.
$ws.Range("E6").Select()
